# Update "want to go" counts (column F) on the "展览" and "全部类型" sheets.
# Each entry is keyed by row number on that sheet -> new value.

$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("展览")
$sheet1Updates = @{
    2  = 1129
    4  = 269
    7  = 543
    8  = 543
    9  = 4995
    13 = 997
    14 = 336
    18 = 3008
    19 = 1854
    23 = 81
    24 = 653
    26 = 317
    28 = 3310
    30 = 2564
    32 = 1605
    33 = 3742
    35 = 903
    36 = 438
    37 = 1172
    39 = 949
    41 = 33
    42 = 912
    43 = 592
    44 = 366
    45 = 377
}
foreach ($row in $sheet1Updates.Keys) {
    $ws1.Range("F$row").Value = $sheet1Updates[$row]
}

$ws4 = $wb.Worksheets.Item("全部类型")
$sheet4Updates = @{
    2  = 1129
    4  = 269
    8  = 543
    9  = 543
    10 = 4995
    15 = 336
    17 = 3008
    19 = 1854
    25 = 81
    27 = 317
    28 = 3310
    32 = 2564
    33 = 1605
    34 = 3742
    37 = 903
    38 = 1172
    40 = 949
    43 = 33
    44 = 912
    45 = 592
    46 = 377
}
foreach ($row in $sheet4Updates.Keys) {
    $ws4.Range("F$row").Value = $sheet4Updates[$row]
}
